$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.326.79"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.642.01"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.03%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "602.34"
$r.ClearFormats()
$ws.Range("E5").Value = "  +1.29%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "146.33"
$r.ClearFormats()
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.108"
$r.ClearFormats()
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "5.61"
$r.ClearFormats()
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.152"
$r.ClearFormats()
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.364"
$r.ClearFormats()
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "27.26"
$r.ClearFormats()
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.113.60"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "63.189.34"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.0000145"
$r.ClearFormats()
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.638.31"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "11.41"
$r.ClearFormats()
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "4.50"
$r.ClearFormats()
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "341.30"
$r.ClearFormats()
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "6.89"
$r.ClearFormats()
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.ClearFormats()
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "5.58"
$r.ClearFormats()
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "66.40"
$r.ClearFormats()
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "1.66"
$r.ClearFormats()
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "8.74"
$r.ClearFormats()
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "554.94"
$r.ClearFormats()
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("B28").Value = "SuiNetwork"
$ws.Range("C28").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "1.52"
$r.ClearFormats()
$ws.Range("E28").Value = "  -4.83%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "0.163"
$r.ClearFormats()
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.ClearFormats()
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "7.82"
$r.ClearFormats()
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "2.01"
$r.ClearFormats()
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.75"
$r.ClearFormats()
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0804"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "5.22"
$r.ClearFormats()
$ws.Range("E35").Value = "  +6.09%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "166.06"
$r.ClearFormats()
$ws.Range("E36").Value = "  -5.28%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.405"
$r.ClearFormats()
$ws.Range("E37").Value = "  +1.07%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.ClearFormats()
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "19.01"
$r.ClearFormats()
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "1.87"
$r.ClearFormats()
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.ClearFormats()
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "168.23"
$r.ClearFormats()
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "22.51"
$r.ClearFormats()
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "3.74"
$r.ClearFormats()
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.0572"
$r.ClearFormats()
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.623"
$r.ClearFormats()
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.0243"
$r.ClearFormats()
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.0960"
$r.ClearFormats()
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "18.71"
$r.ClearFormats()
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "1.76"
$r.ClearFormats()
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "11.23"
$r.ClearFormats()
$ws.Range("E51").Value = "  -1.01%  "
